# Apply the "1 SEPTEMBER ATTENDANCE" edits:
#  - Add a new "Single Class" column-F note to the header row
#  - Record a double-class (F=3) adjustment for the existing roster rows
#  - Append a new student (row 21) with its own attendance + totals

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header note in F10 (new text, formatted like the neighbouring E10 cell) ---
$ws.Range("E10").Copy()
$ws.Range("F10").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("F10").Value = "30 `n (Single `nClass)"

# --- Totals row: one single class recorded in column F ---
$ws.Range("F12").Value = 3

# --- Roster rows 14-19: a "0" is now explicitly recorded in column F ---
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0

# --- Roster row 20 (Imran Khan) attended the single class too ---
$ws.Range("F20").Value = 3

# --- New roster row 21: 22ME-MECH05, single class only ---
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "22ME-MECH05"
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 3
$ws.Range("I21").Formula = "=SUM(D21:H21)"
$ws.Range("J21").Formula = "=(I21/`$I`$12)*100"

# --- Selection / view state left by the editor ---
$ws.Range("D8:H8").Select()
